$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "26.621.20"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.597.80"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.20"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.48"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").Value = "1.821.34"
$ws.Range("D13").Value = "1.595.51"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.78"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "26.604.29"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "0.0₃0738"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.33"
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.01"
$ws.Range("E21").Value = "  +4.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.27"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.33"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.43"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.29"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0512"
$ws.Range("E30").Value = "  +1.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.16"
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("E33").Value = "  +0.89%  "
$ws.Range("D34").Value = "1.282.16"
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("E35").Value = "  -7.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("E38").Value = "  -0.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.839"
$ws.Range("E39").Value = "  +0.35%  "
$ws.Range("E40").Value = "  +21.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.48"
$ws.Range("E41").Value = "  +2.22%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.31"
$ws.Range("E43").Value = "  +1.17%  "
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("D45").Value = "1.733.79"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.22"
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("E48").Value = "  +4.07%  "
$ws.Range("E49").Value = "  +0.76%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.43"
$ws.Range("E51").Value = "  -0.57%  "
